$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new columns right before the old "ausserordentlicheOeffnungszeiten" column (W) ---
# This pushes the old W (and everything to its right) 3 columns to the right.
$ws.Range("W1:Y1").EntireColumn.Insert()

# --- Populate the 4 "new" columns (W,X,Y,Z) in the header row (row 4) and data row (row 5) ---
# Order chosen to reproduce the shared-string table ordering of the authored workbook.
$ws.Range("Y5").Value = "{oeffnungAnWochenenden}"
$ws.Range("Z5").Value = "{uebernachtungMoeglich}"
$ws.Range("Z4").Value = "{uebernachtungMoeglichTitle}"
$ws.Range("Y4").Value = "{oeffnungAnWochenendenTitle}"
$ws.Range("W4").Value = "{oeffnungVorTitle}"
$ws.Range("W5").Value = "{oeffnungVor}"
$ws.Range("X4").Value = "{oeffnungNachTitle}"
$ws.Range("X5").Value = "{oeffnungNach}"

# --- Column widths: W,X share the width of U,V; Y,Z take over the old W's (bestFit) width ---
$ws.Range("W1:X1").ColumnWidth = 19.333333333333336
$ws.Range("Y1:Z1").ColumnWidth = 37.166666666666664

# --- Update view: scroll position + active selection to match the authored state ---
$ws.Range("Z14").Select()
